# Auto-generated edit script: refresh market-board derived columns (H-N)
# on the Pandaemonium_Profits sheets, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 4863.636
$ws.Range("I76").Value = 4388.8887
$ws.Range("J76").Value = 7000
$ws.Range("K76").Value = 4388.8887
$ws.Range("L76").Value = 7000
$ws.Range("M76").Value = -4073.8887
$ws.Range("N76").Value = -7630

# Row 79
$ws.Range("H79").Value = 4863.636
$ws.Range("I79").Value = 4388.8887
$ws.Range("J79").Value = 7000
$ws.Range("K79").Value = 4388.8887
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = -3296.8887
$ws.Range("N79").Value = -9184

# Row 137
$ws.Range("H137").Value = 2113.9822
$ws.Range("I137").Value = 1717.7675
$ws.Range("J137").Value = 3424.5386
$ws.Range("K137").Value = 5153.3025
$ws.Range("L137").Value = 10273.6158
$ws.Range("M137").Value = -2603.3025
$ws.Range("N137").Value = -15373.6158

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6923
$ws.Range("I61").Value = 6079.643
$ws.Range("J61").Value = 7906.9165
$ws.Range("K61").Value = 6079.643
$ws.Range("L61").Value = 7906.9165
$ws.Range("M61").Value = -5867.643
$ws.Range("N61").Value = -8330.916499999999

# Row 74
$ws.Range("H74").Value = 4096.9766
$ws.Range("I74").Value = 1777.9474
$ws.Range("K74").Value = 1777.9474
$ws.Range("M74").Value = -903.9474

# Row 77
$ws.Range("H77").Value = 4096.9766
$ws.Range("I77").Value = 1777.9474
$ws.Range("K77").Value = 8889.737000000001
$ws.Range("M77").Value = -4521.737000000001

# Row 128
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960

# Row 132
$ws.Range("H132").Value = 2386.9473
$ws.Range("I132").Value = 1511.8636
$ws.Range("K132").Value = 4535.5908
$ws.Range("M132").Value = -2005.5908

# Row 136
$ws.Range("H136").Value = 6923
$ws.Range("I136").Value = 6079.643
$ws.Range("J136").Value = 7906.9165
$ws.Range("K136").Value = 18238.929
$ws.Range("L136").Value = 23720.7495
$ws.Range("M136").Value = -15688.929
$ws.Range("N136").Value = -28820.7495

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3219.4783
$ws.Range("I134").Value = 3332.2778
$ws.Range("J134").Value = 2813.4
$ws.Range("K134").Value = 9996.8334
$ws.Range("L134").Value = 8440.200000000001
$ws.Range("M134").Value = -7461.8334
$ws.Range("N134").Value = -13510.2

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 259.44446
$ws.Range("I22").Value = 272.2
$ws.Range("J22").Value = 243.5
$ws.Range("K22").Value = 272.2
$ws.Range("L22").Value = 243.5
$ws.Range("M22").Value = 77.80000000000001
$ws.Range("N22").Value = -943.5

# Row 31
$ws.Range("H31").Value = 2605.2273
$ws.Range("I31").Value = 1790.4706
$ws.Range("J31").Value = 3118.2222
$ws.Range("K31").Value = 1790.4706
$ws.Range("L31").Value = 3118.2222
$ws.Range("M31").Value = -1495.4706
$ws.Range("N31").Value = -3708.2222

# Row 34
$ws.Range("H34").Value = 2605.2273
$ws.Range("I34").Value = 1790.4706
$ws.Range("J34").Value = 3118.2222
$ws.Range("K34").Value = 1790.4706
$ws.Range("L34").Value = 3118.2222
$ws.Range("M34").Value = -1588.4706
$ws.Range("N34").Value = -3522.2222

# Row 58
$ws.Range("H58").Value = 3249469
$ws.Range("I58").Value = 5683189
$ws.Range("J58").Value = 4509.1665
$ws.Range("K58").Value = 5683189
$ws.Range("L58").Value = 4509.1665
$ws.Range("M58").Value = -5682986
$ws.Range("N58").Value = -4915.1665

# Row 132
$ws.Range("H132").Value = 2658.2917
$ws.Range("I132").Value = 2400.303
$ws.Range("J132").Value = 3225.8667
$ws.Range("K132").Value = 7200.909
$ws.Range("L132").Value = 9677.6001
$ws.Range("M132").Value = -4670.909
$ws.Range("N132").Value = -14737.6001

# Row 134
$ws.Range("H134").Value = 2861.9614
$ws.Range("I134").Value = 2683.8948
$ws.Range("J134").Value = 3345.2856
$ws.Range("K134").Value = 8051.6844
$ws.Range("L134").Value = 10035.8568
$ws.Range("M134").Value = -5516.6844
$ws.Range("N134").Value = -15105.8568

# Row 136
$ws.Range("H136").Value = 3249469
$ws.Range("I136").Value = 5683189
$ws.Range("J136").Value = 4509.1665
$ws.Range("K136").Value = 17049567
$ws.Range("L136").Value = 13527.4995
$ws.Range("M136").Value = -17047017
$ws.Range("N136").Value = -18627.4995

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 2560.318
$ws.Range("J46").Value = 2804.375
$ws.Range("L46").Value = 8413.125
$ws.Range("N46").Value = -8595.125

# Row 109
$ws.Range("H109").Value = 2127.0667
$ws.Range("I109").Value = 663.5
$ws.Range("J109").Value = 2659.2727
$ws.Range("K109").Value = 1990.5
$ws.Range("L109").Value = 7977.8181
$ws.Range("M109").Value = -950.5
$ws.Range("N109").Value = -10057.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4706.6665
$ws.Range("I80").Value = 2660
$ws.Range("J80").Value = 5730
$ws.Range("K80").Value = 2660
$ws.Range("L80").Value = 5730
$ws.Range("M80").Value = -1662
$ws.Range("N80").Value = -7726

# Row 83
$ws.Range("H83").Value = 4706.6665
$ws.Range("I83").Value = 2660
$ws.Range("J83").Value = 5730
$ws.Range("K83").Value = 13300
$ws.Range("L83").Value = 28650
$ws.Range("M83").Value = -8308
$ws.Range("N83").Value = -38634

# Row 132
$ws.Range("H132").Value = 3057.4
$ws.Range("I132").Value = 2921.4546
$ws.Range("J132").Value = 3223.5557
$ws.Range("K132").Value = 8764.363799999999
$ws.Range("L132").Value = 9670.667099999999
$ws.Range("M132").Value = -6234.363799999999
$ws.Range("N132").Value = -14730.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 62
$ws.Range("H62").Value = 36874.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 36874.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 36874.5
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -38122.5

# Row 65
$ws.Range("H65").Value = 36874.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 36874.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 110623.5
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -116863.5

# Row 68
$ws.Range("H68").Value = 3851
$ws.Range("J68").Value = 5100
$ws.Range("L68").Value = 5100
$ws.Range("N68").Value = -6598

# Row 71
$ws.Range("H71").Value = 3851
$ws.Range("J71").Value = 5100
$ws.Range("L71").Value = 25500
$ws.Range("N71").Value = -32988

# Row 82
$ws.Range("H82").Value = 1389
$ws.Range("I82").Value = 1116.8334
$ws.Range("J82").Value = 1933.3334
$ws.Range("K82").Value = 1116.8334
$ws.Range("L82").Value = 1933.3334
$ws.Range("M82").Value = -755.8334
$ws.Range("N82").Value = -2655.3334

# Row 85
$ws.Range("H85").Value = 1389
$ws.Range("I85").Value = 1116.8334
$ws.Range("J85").Value = 1933.3334
$ws.Range("K85").Value = 1116.8334
$ws.Range("L85").Value = 1933.3334
$ws.Range("M85").Value = 131.1666
$ws.Range("N85").Value = -4429.3334

# Row 132
$ws.Range("H132").Value = 4528.7754
$ws.Range("I132").Value = 4365.9736
$ws.Range("J132").Value = 5091.1816
$ws.Range("K132").Value = 13097.9208
$ws.Range("L132").Value = 15273.5448
$ws.Range("M132").Value = -10567.9208
$ws.Range("N132").Value = -20333.5448

# Row 136
$ws.Range("H136").Value = 4286.3774
$ws.Range("I136").Value = 2445.7188
$ws.Range("J136").Value = 7091.1904
$ws.Range("K136").Value = 7337.1564
$ws.Range("L136").Value = 21273.5712
$ws.Range("M136").Value = -4787.1564
$ws.Range("N136").Value = -26373.5712

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4066.6667
$ws.Range("J62").Value = 4133.3335
$ws.Range("L62").Value = 4133.3335
$ws.Range("N62").Value = -5381.3335

# Row 65
$ws.Range("H65").Value = 4066.6667
$ws.Range("J65").Value = 4133.3335
$ws.Range("L65").Value = 20666.6675
$ws.Range("N65").Value = -26906.6675

# Row 132
$ws.Range("H132").Value = 1536.1569
$ws.Range("I132").Value = 790.45715
$ws.Range("J132").Value = 3167.375
$ws.Range("K132").Value = 2371.37145
$ws.Range("L132").Value = 9502.125
$ws.Range("M132").Value = 158.6285500000004
$ws.Range("N132").Value = -14562.125

# Row 136
$ws.Range("H136").Value = 6155
$ws.Range("I136").Value = 2504.1667
$ws.Range("J136").Value = 9284.286
$ws.Range("K136").Value = 7512.500100000001
$ws.Range("L136").Value = 27852.858
$ws.Range("M136").Value = -4962.500100000001
$ws.Range("N136").Value = -32952.858
